# Added ability to preview scratch pads before loading them.
# Time Log.xlsx: fill in the (previously blank) row 94 of the log with a
# new "Coding" entry for 2014-10-20, 13:06 -> 14:05 (5 min interruption).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date (stored as serial 41932 = 2014-10-20)
$ws.Range("A94").Value = 41932
# Start / stop time (time-of-day fractions)
$ws.Range("B94").Value = 0.54583333333333328
$ws.Range("C94").Value = 0.58680555555555558
# Interruption, in minutes
$ws.Range("D94").Value = 5
# Delta formula, same shared formula used by the rest of the column
$ws.Range("E94").Formula = "=IF(AND(NOT(ISBLANK(B94)),NOT(ISBLANK(C94))), (C94-B94) * 24 - D94/60, """")"
# Activity/category for this entry
$ws.Range("F94").Value = "Coding"

# Move the active selection down to the next (still blank) row, as if the
# user had just tabbed/entered past the row they filled in.
[void]$ws.Range("A95").Select()
